$d = $word.ActiveDocument

$replacements = @(
    @("18+69=", "63+35="),
    @("13+19=", "88-10="),
    @("37+40=", "48+33="),
    @("28+14=", "67-60="),
    @("6-3=", "26+46="),
    @("72-58=", "41+7="),
    @("75-20=", "43+0="),
    @("97-13=", "22+49="),
    @("86-68=", "76-50="),
    @("71-55=", "79-27="),
    @("98-22=", "54-53="),
    @("88-83=", "36+58="),
    @("35+60=", "89-8="),
    @("21+78=", "83-79="),
    @("9+63=", "30+6="),
    @("13-2=", "57+37="),
    @("17+0=", "20+51="),
    @("57-36=", "63-6="),
    @("31-27=", "24+35="),
    @("15+24=", "80-67="),
    @("17-13=", "93-14="),
    @("59+6=", "11+63="),
    @("93-66=", "4+26="),
    @("34+13=", "25+42="),
    @("98-2=", "95-85="),
    @("11+74=", "54-23="),
    @("69-58=", "93-17="),
    @("60-4=", "71-61="),
    @("76-65=", "30+41="),
    @("8+86=", "40+3="),
    @("15+8=", "30-27="),
    @("47-30=", "10+11="),
    @("58-34=", "69-56="),
    @("73-14=", "90-43="),
    @("53-36=", "35+16="),
    @("81+9=", "1+81="),
    @("82-78=", "30+40="),
    @("89-46=", "9+69="),
    @("46-32=", "70+11="),
    @("58+39=", "33+6="),
    @("15+71=", "77-2="),
    @("71-17=", "61-6="),
    @("47+24=", "24+40="),
    @("76+8=", "99-74="),
    @("67-50=", "8+25="),
    @("70-30=", "84-24="),
    @("64-26=", "76-75="),
    @("19+62=", "18-13="),
    @("11-5=", "6+11="),
    @("35-2=", "8+69="),
    @("4+94=", "90-29="),
    @("53+24=", "8+16="),
    @("91-75=", "79-45="),
    @("91-22=", "96-9="),
    @("90-7=", "55-15="),
    @("26+23=", "31-15="),
    @("62-7=", "5+10="),
    @("30+59=", "47-38="),
    @("65-7=", "88-82="),
    @("32-5=", "20+63="),
    @("29+14=", "9+11="),
    @("73-3=", "85-41="),
    @("72-62=", "21+19="),
    @("81-29=", "69+13="),
    @("37+41=", "21+24="),
    @("72-31=", "11+9="),
    @("68+30=", "42+53="),
    @("57-15=", "19+0="),
    @("12+12=", "45+29="),
    @("97+2=", "34+0="),
    @("25+52=", "99-76="),
    @("67-59=", "87-12="),
    @("18+11=", "95-24="),
    @("68-21=", "87-81="),
    @("0+73=", "96-29="),
    @("19-9=", "56-42="),
    @("72-8=", "74-32="),
    @("40-6=", "4+77="),
    @("20+43=", "26-11="),
    @("95-49=", "87-64="),
    @("50+38=", "27-17="),
    @("41+45=", "90-89="),
    @("71-36=", "35+58="),
    @("59-36=", "79-26="),
    @("98-91=", "33-6="),
    @("77-7=", "39+23="),
    @("37-24=", "39+47="),
    @("37-2=", "2+72="),
    @("95-90=", "99-71="),
    @("53-52=", "80-47="),
    @("18-15=", "10+19="),
    @("81-75=", "12+4="),
    @("45-13=", "50+27="),
    @("44-32=", "56-34="),
    @("66-8=", "31+63="),
    @("3+27=", "22+0="),
    @("72+21=", "3+78="),
    @("54+0=", "10+50="),
    @("28-5=", "58+13="),
    @("76-18=", "73-30=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
